$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Raven Manalastas"

# Test case 1 (row 7)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2016`nmanagement_fee = 2"
$ws.Range("G7").Value = "Attributes are set to input values."

# Test case 2 (row 8)
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "management_fee = ""Invalid fee"""
$ws.Range("G8").Value = "management_fee attribute is set to 2.55"

# Test case 3 (row 9)
$ws.Range("E9").Value = "all inputs must be valid"
$ws.Range("F9").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2016`nmanagement_fee = 2"
$ws.Range("G9").Value = "service_charge = 0.50"

# Test case 4 (row 10)
$ws.Range("E10").Value = "all inputs must be valid"
$ws.Range("F10").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2019`nmanagement_fee = 2"
$ws.Range("G10").Value = "service_charge = 2.50"

# Test case 5 (row 11)
$ws.Range("E11").Value = "all inputs must be valid"
$ws.Range("F11").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2022`nmanagement_fee = 2"
$ws.Range("G11").Value = "service_charge = 2.50"

# Test case 6 (row 12)
$ws.Range("E12").Value = "all inputs must be valid"
$ws.Range("F12").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2016`nmanagement_fee = 2"
$ws.Range("G12").Value = "Account Number: 709 Balance: `$450.00`nDate Created: 2016, 7, 9 Management fee: Waived`nAccount Type: Investment"

# Test case 7 (row 13)
$ws.Range("E13").Value = "all inputs must be valid"
$ws.Range("F13").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2022`nmanagement_fee = 2"
$ws.Range("G13").Value = "Account Number: 709 Balance: `$450.00`nDate Created: 2016, 7, 9 Management fee: `$2.00`nAccount Type: Investment"

# Row heights now taller due to wrapped multi-line content (auto-fit)
for ($r = 7; $r -le 13; $r++) {
    $ws.Rows.Item($r).RowHeight = 84.6
}
$ws.Rows.Item(2).RowHeight = 73.2
for ($r = 14; $r -le 32; $r++) {
    $ws.Rows.Item($r).RowHeight = 31.2
}

# View / selection changes
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 65
$ws.Range("F13").Select() | Out-Null
